$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Conflict")
$ws2.Rows(5).Delete()
